$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 181.4944075
$ws.Range("H2").Value = 362.988815
$ws.Range("I2").Value = 0.2239486468210351
$ws.Range("J2").Value = 0.1654349085470023
$ws.Range("O2").Value = 0.03042272367694056
$ws.Range("P2").Value = 0.04495032978430324
$ws.Range("Q2").Value = 1.342453634141667
$ws.Range("R2").Value = 8.054721804850001
$ws.Range("S2").Value = 0.006813127800061105
$ws.Range("T2").Value = 0.007436353697023802
$ws.Range("G3").Value = 181.4944075
$ws.Range("H3").Value = 362.988815
$ws.Range("I3").Value = 0.2239486468210351
$ws.Range("J3").Value = 0.1654349085470023
$ws.Range("M3").Value = 0.235733
$ws.Range("N3").Value = 0.471466
$ws.Range("O3").Value = 0.9695772763230595
$ws.Range("P3").Value = 0.9550496702156968
$ws.Range("Q3").Value = 42.78422116319751
$ws.Range("R3").Value = 171.13688465279
$ws.Range("S3").Value = 0.217135519020974
$ws.Range("T3").Value = 0.1579985548499785
$ws.Range("I4").Value = 0.07700606288633029
$ws.Range("J4").Value = 0.08532865336765341
$ws.Range("O4").Value = 0.03042272367694056
$ws.Range("P4").Value = 0.04495032978430324
$ws.Range("S4").Value = 0.002342734172639934
$ws.Range("T4").Value = 0.003835551108926518
$ws.Range("I5").Value = 0.07700606288633029
$ws.Range("J5").Value = 0.08532865336765341
$ws.Range("M5").Value = 0.235733
$ws.Range("N5").Value = 0.471466
$ws.Range("O5").Value = 0.9695772763230595
$ws.Range("P5").Value = 0.9550496702156968
$ws.Range("Q5").Value = 14.71160675540367
$ws.Range("R5").Value = 88.269640532422
$ws.Range("S5").Value = 0.07466332871369036
$ws.Range("T5").Value = 0.0814931022587269
$ws.Range("G6").Value = 171.9980316666667
$ws.Range("H6").Value = 515.994095
$ws.Range("I6").Value = 0.2122309275432167
$ws.Range("J6").Value = 0.235168226649403
$ws.Range("O6").Value = 0.03042272367694056
$ws.Range("P6").Value = 0.04495032978430324
$ws.Range("Q6").Value = 1.272212107561111
$ws.Range("R6").Value = 11.44990896805
$ws.Range("S6").Value = 0.006456642864348074
$ws.Range("T6").Value = 0.01057088934268044
$ws.Range("G7").Value = 171.9980316666667
$ws.Range("H7").Value = 515.994095
$ws.Range("I7").Value = 0.2122309275432167
$ws.Range("J7").Value = 0.235168226649403
$ws.Range("M7").Value = 0.235733
$ws.Range("N7").Value = 0.471466
$ws.Range("O7").Value = 0.9695772763230595
$ws.Range("P7").Value = 0.9550496702156968
$ws.Range("Q7").Value = 40.54561199887833
$ws.Range("R7").Value = 243.27367199327
$ws.Range("S7").Value = 0.2057742846788686
$ws.Range("T7").Value = 0.2245973373067226
$ws.Range("G8").Value = 55.64279550000001
$ws.Range("H8").Value = 111.285591
$ws.Range("I8").Value = 0.06865847234198982
$ws.Range("J8").Value = 0.05071925307032974
$ws.Range("O8").Value = 0.03042272367694056
$ws.Range("P8").Value = 0.04495032978430324
$ws.Range("Q8").Value = 0.411571210715
$ws.Range("R8").Value = 2.469427264290001
$ws.Range("S8").Value = 0.002088777732141223
$ws.Range("T8").Value = 0.002279847151924857
$ws.Range("G9").Value = 55.64279550000001
$ws.Range("H9").Value = 111.285591
$ws.Range("I9").Value = 0.06865847234198982
$ws.Range("J9").Value = 0.05071925307032974
$ws.Range("M9").Value = 0.235733
$ws.Range("N9").Value = 0.471466
$ws.Range("O9").Value = 0.9695772763230595
$ws.Range("P9").Value = 0.9550496702156968
$ws.Range("Q9").Value = 13.1168431116015
$ws.Range("R9").Value = 52.467372446406
$ws.Range("S9").Value = 0.0665696946098486
$ws.Range("T9").Value = 0.04843940591840488
$ws.Range("G10").Value = 203.386317
$ws.Range("H10").Value = 610.158951
$ws.Range("I10").Value = 0.250961399315095
$ws.Range("J10").Value = 0.2780845747487284
$ws.Range("O10").Value = 0.03042272367694056
$ws.Range("P10").Value = 0.04495032978430324
$ws.Range("Q10").Value = 1.50438079141
$ws.Range("R10").Value = 13.53942712269
$ws.Range("S10").Value = 0.007634929304941477
$ws.Range("T10").Value = 0.01249999334288307
$ws.Range("G11").Value = 203.386317
$ws.Range("H11").Value = 610.158951
$ws.Range("I11").Value = 0.250961399315095
$ws.Range("J11").Value = 0.2780845747487284
$ws.Range("M11").Value = 0.235733
$ws.Range("N11").Value = 0.471466
$ws.Range("O11").Value = 0.9695772763230595
$ws.Range("P11").Value = 0.9550496702156968
$ws.Range("Q11").Value = 47.944866665361
$ws.Range("R11").Value = 287.669199992166
$ws.Range("S11").Value = 0.2433264700101536
$ws.Range("T11").Value = 0.2655845814058454
$ws.Range("G12").Value = 135.4992116666667
$ws.Range("H12").Value = 406.497635
$ws.Range("I12").Value = 0.167194491092333
$ws.Range("J12").Value = 0.1852643836168829
$ws.Range("O12").Value = 0.03042272367694056
$ws.Range("P12").Value = 0.04495032978430324
$ws.Range("Q12").Value = 1.002242502294445
$ws.Range("R12").Value = 9.020182520650001
$ws.Range("S12").Value = 0.005086511802808748
$ws.Range("T12").Value = 0.008327695140864552
$ws.Range("G13").Value = 135.4992116666667
$ws.Range("H13").Value = 406.497635
$ws.Range("I13").Value = 0.167194491092333
$ws.Range("J13").Value = 0.1852643836168829
$ws.Range("M13").Value = 0.235733
$ws.Range("N13").Value = 0.471466
$ws.Range("O13").Value = 0.9695772763230595
$ws.Range("P13").Value = 0.9550496702156968
$ws.Range("Q13").Value = 31.94163566381833
$ws.Range("R13").Value = 191.64981398291
$ws.Range("S13").Value = 0.1621079792895243
$ws.Range("T13").Value = 0.1769366884760183
